# Apply cryptos list update (Sat Sep 16 22:26:40 UTC 2023 GitHub Actions run).
# Source: unified diff of canonical OOXML for cryptos.xlsx (sheet1), rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that LOOK numeric (e.g. "216.53", "0.0520").
# Assigning such a string straight to .Value lets Excel auto-coerce it to a
# real number (losing the original text formatting/precision, e.g. trailing
# zeros). The source file stores these as plain text, so we force text entry
# by switching the cell to the "@" (Text) number format before the write, then
# clear the format again immediately after so the cell keeps its original
# (default) style - only the stored value/type changes.
function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue "D2" "26.801.94"
# Row 3
Set-TextValue "D3" "1.649.32"
$ws.Range("E3").Value = "  +0.75%  "
# Row 4
$ws.Range("E4").Value = "  +0.61%  "
# Row 5
Set-TextValue "D5" "216.53"
$ws.Range("E5").Value = "  +1.22%  "
# Row 6
$ws.Range("E6").Value = "  +0.62%  "
# Row 7
$ws.Range("E7").Value = "  +0.55%  "
# Row 8
$ws.Range("E8").Value = "  +0.57%  "
# Row 9
$ws.Range("E9").Value = "  +0.11%  "
# Row 10
Set-TextValue "D10" "19.25"
$ws.Range("E10").Value = "  +1.66%  "
# Row 11
Set-TextValue "D11" "0.0844"
$ws.Range("E11").Value = "  +0.09%  "
# Row 12
Set-TextValue "D12" "1.879.16"
$ws.Range("E12").Value = "  +0.79%  "
# Row 13
Set-TextValue "D13" "1.665.65"
$ws.Range("E13").Value = "  +2.48%  "
# Row 14
$ws.Range("E14").Value = "  +0.94%  "
# Row 15
$ws.Range("E15").Value = "  +1.07%  "
# Row 16
Set-TextValue "D16" "65.50"
$ws.Range("E16").Value = "  +0.17%  "
# Row 17
Set-TextValue "D17" "26.799.69"
# Row 18
Set-TextValue "D18" "0.0₃0744"
$ws.Range("E18").Value = "  +0.08%  "
# Row 19
Set-TextValue "D19" "218.40"
$ws.Range("E19").Value = "  +0.76%  "
# Row 20
$ws.Range("E20").Value = "  +0.56%  "
# Row 21
$ws.Range("E21").Value = "  +11.91%  "
# Row 22
Set-TextValue "D22" "4.37"
$ws.Range("E22").Value = "  +0.95%  "
# Row 23
$ws.Range("E23").Value = "  +0.72%  "
# Row 24
Set-TextValue "D24" "9.46"
$ws.Range("E24").Value = "  +0.84%  "
# Row 25
Set-TextValue "D25" "145.91"
$ws.Range("E25").Value = "  -0.97%  "
# Row 26
$ws.Range("E26").Value = "  +0.58%  "
# Row 27
$ws.Range("E27").Value = "  -0.01%  "
# Row 28
$ws.Range("E28").Value = "  +4.31%  "
# Row 29
Set-TextValue "D29" "15.83"
$ws.Range("E29").Value = "  +0.98%  "
# Row 30
Set-TextValue "D30" "0.0519"
$ws.Range("E30").Value = "  +0.61%  "
# Row 31
$ws.Range("E31").Value = "  +1.18%  "
# Row 32
$ws.Range("E32").Value = "  -0.35%  "
# Row 33
$ws.Range("E33").Value = "  +0.65%  "
# Row 34
Set-TextValue "D34" "1.283.91"
$ws.Range("E34").Value = "  +1.01%  "
# Row 36
$ws.Range("E36").Value = "  +2.40%  "
# Row 37
$ws.Range("E37").Value = "  +1.69%  "
# Row 38
Set-TextValue "D38" "0.541"
$ws.Range("E38").Value = "  +5.53%  "
# Row 39
$ws.Range("E39").Value = "  +3.85%  "
# Row 40
$ws.Range("E40").Value = "  +0.61%  "
# Row 41
$ws.Range("E41").Value = "  +2.24%  "
# Row 42
$ws.Range("E42").Value = "  -1.00%  "
# Row 43
$ws.Range("E43").Value = "  +1.67%  "
# Row 44
Set-TextValue "D44" "1.790.67"
$ws.Range("E44").Value = "  +0.96%  "
# Row 45
Set-TextValue "D45" "92.04"
$ws.Range("E45").Value = "  -1.64%  "
# Row 46
Set-TextValue "D46" "59.84"
$ws.Range("E46").Value = "  +8.25%  "
# Row 47
Set-TextValue "D47" "1.62"
$ws.Range("E47").Value = "  +0.95%  "
# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0516"
$ws.Range("E48").Value = "  +0.72%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "7.81"
$ws.Range("E49").Value = "  +2.45%  "
# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0981"
$ws.Range("E50").Value = "  +2.02%  "
# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.409"
$ws.Range("E51").Value = "  +0.31%  "
